$wb = $excel.ActiveWorkbook

# --- Sheet2: add header row ---
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("A1").Value = "Header"
$ws2.Range("B1").Value = "X Min"
$ws2.Range("C1").Value = "X Max"
$ws2.Range("D1").Value = "Y Min"
$ws2.Range("E1").Value = "Y Max"
$ws2.Range("E1").Select()

# --- Sheet1: add Range column (G) ---
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("G1").Value = "Range"
$ws1.Range("G2").Value = "9:20"
$ws1.Range("G3").Value = "2:5"
$ws1.Range("G4").Value = "3:6"
$ws1.Range("G6").Value = "1:4"

# Apply text number format to G2:G6 (numFmtId 49 = "@" text format)
$ws1.Range("G2:G6").NumberFormat = "@"

# Set page orientation to portrait (adds pageSetup element)
$ws1.PageSetup.Orientation = 1

$ws1.Range("G5").Select()
